$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sets a cell to a text value without letting Excel auto-coerce numeric-looking
# strings (prices like "1.000" or "0.00001000") into actual numbers. Using an
# apostrophe text-prefix forces text entry; the Style is then reset back to
# "Normal" so the quote-prefix formatting does not leave a stray style behind.
function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextCell "D2" "29.282.96"
Set-TextCell "E2" "  +0.43%  "

# Row 3
Set-TextCell "D3" "1.896.98"
Set-TextCell "E3" "  -0.60%  "

# Row 4
Set-TextCell "D4" "1.000"
Set-TextCell "E4" "  -0.46%  "

# Row 5
Set-TextCell "D5" "325.75"
Set-TextCell "E5" "  -0.73%  "

# Row 6
Set-TextCell "D6" "0.9984"
Set-TextCell "E6" "  -0.55%  "

# Row 7
Set-TextCell "D7" "0.4633"
Set-TextCell "E7" "  +0.02%  "

# Row 8
Set-TextCell "D8" "0.3913"
Set-TextCell "E8" "  -0.48%  "

# Row 9
Set-TextCell "D9" "0.07887"
Set-TextCell "E9" "  -1.24%  "

# Row 10
Set-TextCell "D10" "0.9909"
Set-TextCell "E10" "  -1.59%  "

# Row 11
Set-TextCell "D11" "22.00"
Set-TextCell "E11" "  -1.75%  "

# Row 12
Set-TextCell "D12" "1.899.05"
Set-TextCell "E12" "  -11.61%  "

# Row 13
Set-TextCell "D13" "7.085"
Set-TextCell "E13" "  -1.57%  "

# Row 14
Set-TextCell "D14" "5.759"
Set-TextCell "E14" "  -0.28%  "

# Row 15
Set-TextCell "D15" "0.06982"
Set-TextCell "E15" "  +0.02%  "

# Row 16
Set-TextCell "D16" "88.78"
Set-TextCell "E16" "  -0.26%  "

# Row 17
Set-TextCell "D17" "0.9992"
Set-TextCell "E17" "  -0.57%  "

# Row 18
Set-TextCell "D18" "0.00001000"
Set-TextCell "E18" "  -1.16%  "

# Row 19
Set-TextCell "D19" "17.13"
Set-TextCell "E19" "  -0.90%  "

# Row 20
Set-TextCell "D20" "0.9971"
Set-TextCell "E20" "  -0.64%  "

# Row 21
Set-TextCell "D21" "29.268.49"
Set-TextCell "E21" "  +0.21%  "

# Row 22
Set-TextCell "D22" "5.298"
Set-TextCell "E22" "  -1.84%  "

# Row 23
Set-TextCell "D23" "11.07"
Set-TextCell "E23" "  -0.25%  "

# Row 24
Set-TextCell "D24" "2.099"
Set-TextCell "E24" "  +2.33%  "

# Row 25
Set-TextCell "D25" "155.89"
Set-TextCell "E25" "  -0.73%  "

# Row 26
Set-TextCell "D26" "19.43"
Set-TextCell "E26" "  -0.95%  "

# Row 27
Set-TextCell "D27" "6.016"
Set-TextCell "E27" "  +2.03%  "

# Row 28
Set-TextCell "D28" "118.51"
Set-TextCell "E28" "  -1.15%  "

# Row 29
Set-TextCell "D29" "1.920"
Set-TextCell "E29" "  -4.62%  "

# Row 30
Set-TextCell "D30" "0.09369"
Set-TextCell "E30" "  +0.16%  "

# Row 31
Set-TextCell "D31" "0.9081"
Set-TextCell "E31" "  -2.31%  "

# Row 32
Set-TextCell "D32" "5.299"
Set-TextCell "E32" "  -1.36%  "

# Row 33
Set-TextCell "D33" "1.330"
Set-TextCell "E33" "  -1.33%  "

# Row 34
Set-TextCell "D34" "3.224"
Set-TextCell "E34" "  -1.72%  "

# Row 35
Set-TextCell "D35" "0.05800"
Set-TextCell "E35" "  -0.81%  "

# Row 36
Set-TextCell "D36" "1.177"
Set-TextCell "E36" "  +1.57%  "

# Row 37
Set-TextCell "D37" "0.02089"
Set-TextCell "E37" "  -0.75%  "

# Row 38
Set-TextCell "B38" "FraxShare"
Set-TextCell "C38" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D38" "7.794"
Set-TextCell "E38" "  -3.28%  "

# Row 39
Set-TextCell "B39" "Frax"
Set-TextCell "C39" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D39" "0.9969"
Set-TextCell "E39" "  -0.70%  "

# Row 40
Set-TextCell "D40" "0.5711"
Set-TextCell "E40" "  -1.28%  "

# Row 41
Set-TextCell "D41" "0.1785"
Set-TextCell "E41" "  -1.61%  "

# Row 42
Set-TextCell "D42" "9.755"
Set-TextCell "E42" "  -2.56%  "

# Row 43
Set-TextCell "D43" "11.94"
Set-TextCell "E43" "  -0.79%  "

# Row 44
Set-TextCell "D44" "2.241"
Set-TextCell "E44" "  -0.74%  "

# Row 45
Set-TextCell "D45" "0.5362"
Set-TextCell "E45" "  -1.26%  "

# Row 46
Set-TextCell "D46" "0.07053"
Set-TextCell "E46" "  -1.18%  "

# Row 47
Set-TextCell "D47" "1.860"
Set-TextCell "E47" "  -1.53%  "

# Row 48
Set-TextCell "B48" "MXToken"
Set-TextCell "C48" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D48" "2.557"
Set-TextCell "E48" "  +1.48%  "

# Row 49
Set-TextCell "B49" "Quant"
Set-TextCell "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D49" "113.10"
Set-TextCell "E49" "  +0.53%  "

# Row 50
Set-TextCell "D50" "1.071"
Set-TextCell "E50" "  -4.43%  "

# Row 51
Set-TextCell "B51" "Aave"
Set-TextCell "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D51" "71.52"
Set-TextCell "E51" "  -0.73%  "

